$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The HED tag for the "PerturbLeft" row was replaced with a shorter tag.
$ws.Range("D2").Value = "Attribute/Sensory/Bisual"

# With the shorter text the wrapped row no longer needs the old custom
# height, so let Excel recompute/auto-fit it instead of leaving the old
# pinned height in place.
$ws.Rows.Item(2).AutoFit()

# The active/selected cell at save time moved to D3.
$ws.Range("D3").Select()
